$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Recomputed NATMI ligand/receptor statistics ("Natmi following Dr Hou advice").
# The ligand/receptor-expressing cell counts (columns E and K) went from 1 to 3,
# which cascades into the dependent specificity/expression-weight columns
# (G,H,I,J,M,N,O,P,Q,R,S,T) for every data row (2-9). Columns F and L are unchanged.

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 1.558564
$ws.Cells.Item(2, 8).Value = 4.675692
$ws.Cells.Item(2, 9).Value = 0.005692101168584756
$ws.Cells.Item(2, 10).Value = 0.005692101168584756
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 45.931316
$ws.Cells.Item(2, 14).Value = 137.793948
$ws.Cells.Item(2, 15).Value = 0.9874217014725413
$ws.Cells.Item(2, 16).Value = 0.9874217014725412
$ws.Cells.Item(2, 17).Value = 71.586895590224
$ws.Cells.Item(2, 18).Value = 644.2820603120159
$ws.Cells.Item(2, 19).Value = 0.0056205042208378
$ws.Cells.Item(2, 20).Value = 0.0056205042208378
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 1.558564
$ws.Cells.Item(3, 8).Value = 4.675692
$ws.Cells.Item(3, 9).Value = 0.005692101168584756
$ws.Cells.Item(3, 10).Value = 0.005692101168584756
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 0.5850973333333334
$ws.Cells.Item(3, 14).Value = 1.755292
$ws.Cells.Item(3, 15).Value = 0.01257829852745884
$ws.Cells.Item(3, 16).Value = 0.01257829852745884
$ws.Cells.Item(3, 17).Value = 0.9119116402293332
$ws.Cells.Item(3, 18).Value = 8.207204762064
$ws.Cells.Item(3, 19).Value = 0.00007159694774695638
$ws.Cells.Item(3, 20).Value = 0.00007159694774695638
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 259.5505726666667
$ws.Cells.Item(4, 8).Value = 778.6517180000001
$ws.Cells.Item(4, 9).Value = 0.9479162344201305
$ws.Cells.Item(4, 10).Value = 0.9479162344201304
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 45.931316
$ws.Cells.Item(4, 14).Value = 137.793948
$ws.Cells.Item(4, 15).Value = 0.9874217014725413
$ws.Cells.Item(4, 16).Value = 0.9874217014725412
$ws.Cells.Item(4, 17).Value = 11921.49937113363
$ws.Cells.Item(4, 18).Value = 107293.4943402027
$ws.Cells.Item(4, 19).Value = 0.9359930610445696
$ws.Cells.Item(4, 20).Value = 0.9359930610445694
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 259.5505726666667
$ws.Cells.Item(5, 8).Value = 778.6517180000001
$ws.Cells.Item(5, 9).Value = 0.9479162344201305
$ws.Cells.Item(5, 10).Value = 0.9479162344201304
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 0.5850973333333334
$ws.Cells.Item(5, 14).Value = 1.755292
$ws.Cells.Item(5, 15).Value = 0.01257829852745884
$ws.Cells.Item(5, 16).Value = 0.01257829852745884
$ws.Cells.Item(5, 17).Value = 151.8623479324063
$ws.Cells.Item(5, 18).Value = 1366.761131391656
$ws.Cells.Item(5, 19).Value = 0.01192317337556105
$ws.Cells.Item(5, 20).Value = 0.01192317337556105
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 12.18925266666667
$ws.Cells.Item(6, 8).Value = 36.567758
$ws.Cells.Item(6, 9).Value = 0.04451691386950307
$ws.Cells.Item(6, 10).Value = 0.04451691386950307
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 45.931316
$ws.Cells.Item(6, 14).Value = 137.793948
$ws.Cells.Item(6, 15).Value = 0.9874217014725413
$ws.Cells.Item(6, 16).Value = 0.9874217014725412
$ws.Cells.Item(6, 17).Value = 559.8684160365093
$ws.Cells.Item(6, 18).Value = 5038.815744328584
$ws.Cells.Item(6, 19).Value = 0.04395696683733129
$ws.Cells.Item(6, 20).Value = 0.04395696683733129
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 12.18925266666667
$ws.Cells.Item(7, 8).Value = 36.567758
$ws.Cells.Item(7, 9).Value = 0.04451691386950307
$ws.Cells.Item(7, 10).Value = 0.04451691386950307
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 0.5850973333333334
$ws.Cells.Item(7, 14).Value = 1.755292
$ws.Cells.Item(7, 15).Value = 0.01257829852745884
$ws.Cells.Item(7, 16).Value = 0.01257829852745884
$ws.Cells.Item(7, 17).Value = 7.131899230592889
$ws.Cells.Item(7, 18).Value = 64.18709307533599
$ws.Cells.Item(7, 19).Value = 0.0005599470321717825
$ws.Cells.Item(7, 20).Value = 0.0005599470321717825
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 0.5133286666666667
$ws.Cells.Item(8, 8).Value = 1.539986
$ws.Cells.Item(8, 9).Value = 0.001874750541781658
$ws.Cells.Item(8, 10).Value = 0.001874750541781658
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 45.931316
$ws.Cells.Item(8, 14).Value = 137.793948
$ws.Cells.Item(8, 15).Value = 0.9874217014725413
$ws.Cells.Item(8, 16).Value = 0.9874217014725412
$ws.Cells.Item(8, 17).Value = 23.57786120052533
$ws.Cells.Item(8, 18).Value = 212.200750804728
$ws.Cells.Item(8, 19).Value = 0.001851169369802613
$ws.Cells.Item(8, 20).Value = 0.001851169369802613
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 0.5133286666666667
$ws.Cells.Item(9, 8).Value = 1.539986
$ws.Cells.Item(9, 9).Value = 0.001874750541781658
$ws.Cells.Item(9, 10).Value = 0.001874750541781658
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 0.5850973333333334
$ws.Cells.Item(9, 14).Value = 1.755292
$ws.Cells.Item(9, 15).Value = 0.01257829852745884
$ws.Cells.Item(9, 16).Value = 0.01257829852745884
$ws.Cells.Item(9, 17).Value = 0.3003472339902222
$ws.Cells.Item(9, 18).Value = 2.703125105912
$ws.Cells.Item(9, 19).Value = 0.00002358117197904489
$ws.Cells.Item(9, 20).Value = 0.00002358117197904489
